$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 192
$wsExhibit.Range("F4").Value = 2396
$wsExhibit.Range("F5").Value = 30

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 192
$wsAll.Range("F6").Value = 2396
$wsAll.Range("F7").Value = 30
